$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D7 with new ID value
$ws.Range("D7").Value = "1245789663"

# Row 8 now holds what used to be row 9's data
$ws.Range("A8").Value = "Lic. Albertino Jesus Heredia Paez"
$ws.Range("B8").Value = "Alexander Francisco Tibanta Miranda"
$ws.Range("C8").Value = "1728220441001"
$ws.Range("D8").Value = "1728220441"
$ws.Range("E8").Value = "Dr. Christian Santiago Izurieta Cruz"

# Row 9 now holds what used to be row 10's data
$ws.Range("A9").Value = "Lic. Alexander Javier Miranda Granero"
$ws.Range("B9").Value = "Peter Patricio Tene Ojeda"
$ws.Range("C9").Value = "174582556"
$ws.Range("D9").Value = "174582556001"
$ws.Range("E9").Value = "Dr. Christian Santiago Izurieta Cruz"

# Remove the old row 10 entirely (shift rows up)
$ws.Range("A10:E10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
